{"js": "// Update the date line and every two-digit \u00d7 two-digit multiplication\n// prompt in the practice-sheet table to the newly generated values.\nconst replacements = [\n  [\"2024-03-11 Monday\", \"2024-03-12 Tuesday\"],\n  [\"24\\u00D779=\", \"35\\u00D760=\"],\n  [\"70\\u00D711=\", \"15\\u00D734=\"],\n  [\"26\\u00D799=\", \"66\\u00D729=\"],\n  [\"11\\u00D747=\", \"36\\u00D786=\"],\n  [\"89\\u00D756=\", \"94\\u00D758=\"],\n  [\"58\\u00D768=\", \"97\\u00D718=\"],\n  [\"79\\u00D789=\", \"50\\u00D750=\"],\n  [\"99\\u00D740=\", \"28\\u00D751=\"],\n  [\"67\\u00D732=\", \"66\\u00D758=\"],\n  [\"33\\u00D767=\", \"53\\u00D750=\"],\n  [\"67\\u00D782=\", \"44\\u00D723=\"],\n  [\"45\\u00D766=\", \"26\\u00D745=\"],\n  [\"95\\u00D794=\", \"37\\u00D730=\"],\n  [\"86\\u00D725=\", \"95\\u00D719=\"],\n  [\"71\\u00D788=\", \"37\\u00D714=\"],\n  [\"59\\u00D739=\", \"94\\u00D762=\"],\n  [\"92\\u00D794=\", \"33\\u00D777=\"],\n  [\"35\\u00D763=\", \"19\\u00D785=\"],\n  [\"45\\u00D762=\", \"30\\u00D712=\"],\n  [\"43\\u00D714=\", \"33\\u00D727=\"],\n  [\"50\\u00D739=\", \"26\\u00D731=\"],\n  [\"53\\u00D727=\", \"20\\u00D770=\"],\n  [\"40\\u00D719=\", \"59\\u00D778=\"],\n  [\"54\\u00D776=\", \"33\\u00D717=\"],\n  [\"45\\u00D741=\", \"92\\u00D719=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit x two-digit multiplication\n# prompt in the practice-sheet table to the newly generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-11 Monday\", \"2024-03-12 Tuesday\"),\n    @(\"24\u00d779=\", \"35\u00d760=\"),\n    @(\"70\u00d711=\", \"15\u00d734=\"),\n    @(\"26\u00d799=\", \"66\u00d729=\"),\n    @(\"11\u00d747=\", \"36\u00d786=\"),\n    @(\"89\u00d756=\", \"94\u00d758=\"),\n    @(\"58\u00d768=\", \"97\u00d718=\"),\n    @(\"79\u00d789=\", \"50\u00d750=\"),\n    @(\"99\u00d740=\", \"28\u00d751=\"),\n    @(\"67\u00d732=\", \"66\u00d758=\"),\n    @(\"33\u00d767=\", \"53\u00d750=\"),\n    @(\"67\u00d782=\", \"44\u00d723=\"),\n    @(\"45\u00d766=\", \"26\u00d745=\"),\n    @(\"95\u00d794=\", \"37\u00d730=\"),\n    @(\"86\u00d725=\", \"95\u00d719=\"),\n    @(\"71\u00d788=\", \"37\u00d714=\"),\n    @(\"59\u00d739=\", \"94\u00d762=\"),\n    @(\"92\u00d794=\", \"33\u00d777=\"),\n    @(\"35\u00d763=\", \"19\u00d785=\"),\n    @(\"45\u00d762=\", \"30\u00d712=\"),\n    @(\"43\u00d714=\", \"33\u00d727=\"),\n    @(\"50\u00d739=\", \"26\u00d731=\"),\n    @(\"53\u00d727=\", \"20\u00d770=\"),\n    @(\"40\u00d719=\", \"59\u00d778=\"),\n    @(\"54\u00d776=\", \"33\u00d717=\"),\n    @(\"45\u00d741=\", \"92\u00d719=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
